# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on several leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 5324.846
$ws.Range("I28").Value = 332.5
$ws.Range("J28").Value = 16557.625
$ws.Range("K28").Value = 332.5
$ws.Range("L28").Value = 16557.625
$ws.Range("M28").Value = 152.5
$ws.Range("N28").Value = -17527.625
$ws.Range("H62").Value = 3860
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3860
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = 3860
$ws.Range("N62").Value = -5108
$ws.Range("H65").Value = 3860
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3860
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = 19300
$ws.Range("N65").Value = -25540
$ws.Range("H107").Value = 942.6667
$ws.Range("I107").Value = 869.1818
$ws.Range("J107").Value = 1751
$ws.Range("K107").Value = 869.1818
$ws.Range("L107").Value = 1751
$ws.Range("M107").Value = 1050.8182
$ws.Range("N107").Value = -5591
$ws.Range("H108").Value = 29900
$ws.Range("J108").Value = 29900
$ws.Range("L108").Value = 29900
$ws.Range("N108").Value = -37580
$ws.Range("H132").Value = 5335.5713
$ws.Range("I132").Value = 6049.0557
$ws.Range("J132").Value = 4051.3
$ws.Range("K132").Value = 18147.1671
$ws.Range("L132").Value = 12153.9
$ws.Range("M132").Value = -15617.1671
$ws.Range("N132").Value = -17213.9
$ws.Range("H135").Value = 1878.5264
$ws.Range("I135").Value = 1878.5264
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 16906.7376
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -14371.7376
$ws.Range("H137").Value = 4605.7144
$ws.Range("I137").Value = 4906.154
$ws.Range("J137").Value = 700
$ws.Range("K137").Value = 14718.462
$ws.Range("L137").Value = 2100
$ws.Range("M137").Value = -12168.462
$ws.Range("N137").Value = -7200
$ws.Range("H138").Value = 2009.5778
$ws.Range("I138").Value = 1798.92
$ws.Range("J138").Value = 2090.6
$ws.Range("K138").Value = 5396.76
$ws.Range("L138").Value = 6271.799999999999
$ws.Range("M138").Value = -256.7600000000002
$ws.Range("N138").Value = -16551.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1396.0769
$ws.Range("I2").Value = 990
$ws.Range("J2").Value = 1649.875
$ws.Range("K2").Value = 990
$ws.Range("L2").Value = 1649.875
$ws.Range("M2").Value = -877
$ws.Range("N2").Value = -1875.875
$ws.Range("H34").Value = 30009.334
$ws.Range("J34").Value = 30009.334
$ws.Range("L34").Value = 30009.334
$ws.Range("N34").Value = -30551.334
$ws.Range("H45").Value = 4289.2856
$ws.Range("I45").Value = 4737.5
$ws.Range("J45").Value = 4110
$ws.Range("K45").Value = 4737.5
$ws.Range("L45").Value = 4110
$ws.Range("M45").Value = -4360.5
$ws.Range("N45").Value = -4864
$ws.Range("H96").Value = 49547
$ws.Range("J96").Value = 49547
$ws.Range("L96").Value = 49547
$ws.Range("N96").Value = -55039
$ws.Range("H116").Value = 1396.0769
$ws.Range("I116").Value = 990
$ws.Range("J116").Value = 1649.875
$ws.Range("K116").Value = 990
$ws.Range("L116").Value = 1649.875
$ws.Range("M116").Value = 1304
$ws.Range("N116").Value = -6237.875
$ws.Range("H132").Value = 5164.3784
$ws.Range("I132").Value = 5038.0435
$ws.Range("J132").Value = 5371.9287
$ws.Range("K132").Value = 15114.1305
$ws.Range("L132").Value = 16115.7861
$ws.Range("M132").Value = -12584.1305
$ws.Range("N132").Value = -21175.7861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1396.0769
$ws.Range("I3").Value = 990
$ws.Range("J3").Value = 1649.875
$ws.Range("K3").Value = 990
$ws.Range("L3").Value = 1649.875
$ws.Range("M3").Value = -876
$ws.Range("N3").Value = -1877.875
$ws.Range("H38").Value = 9000
$ws.Range("J38").Value = 9000
$ws.Range("L38").Value = 9000
$ws.Range("N38").Value = -9832
$ws.Range("H105").Value = 6946319
$ws.Range("I105").Value = 7814234
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 7814234
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -7812487
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3178.0417
$ws.Range("I31").Value = 960.32434
$ws.Range("J31").Value = 10637.637
$ws.Range("K31").Value = 960.32434
$ws.Range("L31").Value = 10637.637
$ws.Range("M31").Value = -665.32434
$ws.Range("N31").Value = -11227.637
$ws.Range("H34").Value = 3178.0417
$ws.Range("I34").Value = 960.32434
$ws.Range("J34").Value = 10637.637
$ws.Range("K34").Value = 960.32434
$ws.Range("L34").Value = 10637.637
$ws.Range("M34").Value = -758.32434
$ws.Range("N34").Value = -11041.637

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 14492.714
$ws.Range("I33").Value = 16883.166
$ws.Range("K33").Value = 101298.996
$ws.Range("M33").Value = -101015.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2544.2222
$ws.Range("I97").Value = 2682
$ws.Range("J97").Value = 2372
$ws.Range("K97").Value = 2682
$ws.Range("L97").Value = 2372
$ws.Range("M97").Value = -2186
$ws.Range("N97").Value = -3364

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3425.7856
$ws.Range("I61").Value = 3432.2632
$ws.Range("J61").Value = 3412.111
$ws.Range("K61").Value = 3432.2632
$ws.Range("L61").Value = 3412.111
$ws.Range("M61").Value = -3230.2632
$ws.Range("N61").Value = -3816.111
$ws.Range("H113").Value = 3425.7856
$ws.Range("I113").Value = 3432.2632
$ws.Range("J113").Value = 3412.111
$ws.Range("K113").Value = 3432.2632
$ws.Range("L113").Value = 3412.111
$ws.Range("M113").Value = -1262.2632
$ws.Range("N113").Value = -7752.111
$ws.Range("H136").Value = 6412220
$ws.Range("I136").Value = 2256.4285
$ws.Range("K136").Value = 6769.2855
$ws.Range("M136").Value = -4219.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H126").Value = 1047.4762
$ws.Range("I126").Value = 826.13336
$ws.Range("J126").Value = 1600.8334
$ws.Range("K126").Value = 2478.40008
$ws.Range("L126").Value = 4802.5002
$ws.Range("M126").Value = -8.400080000000344
$ws.Range("N126").Value = -9742.5002
$ws.Range("H136").Value = 3141.5173
$ws.Range("I136").Value = 2764.5715
$ws.Range("J136").Value = 3493.3333
$ws.Range("K136").Value = 8293.7145
$ws.Range("L136").Value = 10479.9999
$ws.Range("M136").Value = -5743.7145
$ws.Range("N136").Value = -15579.9999
